# Auto-generated: apply row-content swaps/permutation for South Korea K League 2 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Cells.Item(11, 2).Value = 6206186
$ws.Cells.Item(11, 6).Value = 'Bucheon'
$ws.Cells.Item(11, 7).Value = 'Cheonan City'
$ws.Cells.Item(11, 8).Value = 1
$ws.Cells.Item(11, 9).Value = 1
$ws.Cells.Item(11, 10).Value = 'D'
$ws.Cells.Item(11, 11).Value = 1.5
$ws.Cells.Item(11, 12).Value = 4
$ws.Cells.Item(11, 13).Value = 5.25
$ws.Cells.Item(11, 14).Value = 1.5
$ws.Cells.Item(11, 15).Value = 4
$ws.Cells.Item(11, 16).Value = 5
$ws.Cells.Item(11, 17).Value = -1
$ws.Cells.Item(11, 18).Value = 1.875
$ws.Cells.Item(11, 19).Value = 1.925
$ws.Cells.Item(11, 20).Value = 2.5
$ws.Cells.Item(11, 21).Value = 1.875
$ws.Cells.Item(11, 22).Value = 1.925
$ws.Cells.Item(11, 23).Value = -1
$ws.Cells.Item(11, 24).Value = 3
$ws.Cells.Item(11, 25).Value = -1
$ws.Cells.Item(11, 26).Value = -1
$ws.Cells.Item(11, 27).Value = 0.925
$ws.Cells.Item(11, 28).Value = -1
$ws.Cells.Item(11, 29).Value = 0.925

# Row 12
$ws.Cells.Item(12, 2).Value = 6204719
$ws.Cells.Item(12, 6).Value = 'Chungbuk Cheongju'
$ws.Cells.Item(12, 7).Value = 'Jeonnam Dragons'
$ws.Cells.Item(12, 8).Value = 3
$ws.Cells.Item(12, 9).Value = 1
$ws.Cells.Item(12, 10).Value = 'H'
$ws.Cells.Item(12, 11).Value = 3.5
$ws.Cells.Item(12, 12).Value = 3.3
$ws.Cells.Item(12, 13).Value = 1.95
$ws.Cells.Item(12, 14).Value = 3.1
$ws.Cells.Item(12, 15).Value = 3.1
$ws.Cells.Item(12, 16).Value = 2.2
$ws.Cells.Item(12, 17).Value = 0.25
$ws.Cells.Item(12, 18).Value = 1.825
$ws.Cells.Item(12, 19).Value = 1.975
$ws.Cells.Item(12, 20).Value = 2.25
$ws.Cells.Item(12, 21).Value = 1.825
$ws.Cells.Item(12, 22).Value = 1.975
$ws.Cells.Item(12, 23).Value = 2.1
$ws.Cells.Item(12, 24).Value = -1
$ws.Cells.Item(12, 25).Value = -1
$ws.Cells.Item(12, 26).Value = 0.825
$ws.Cells.Item(12, 27).Value = -1
$ws.Cells.Item(12, 28).Value = 0.825
$ws.Cells.Item(12, 29).Value = -1

# Row 24
$ws.Cells.Item(24, 2).Value = 6206197
$ws.Cells.Item(24, 6).Value = 'Chungbuk Cheongju'
$ws.Cells.Item(24, 7).Value = 'Seongnam FC'
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(24, 9).Value = 0
$ws.Cells.Item(24, 10).Value = 'D'
$ws.Cells.Item(24, 11).Value = 2.7
$ws.Cells.Item(24, 12).Value = 3.1
$ws.Cells.Item(24, 13).Value = 2.5
$ws.Cells.Item(24, 14).Value = 2.7
$ws.Cells.Item(24, 15).Value = 3
$ws.Cells.Item(24, 16).Value = 2.55
$ws.Cells.Item(24, 17).Value = 0
$ws.Cells.Item(24, 18).Value = 1.925
$ws.Cells.Item(24, 19).Value = 1.875
$ws.Cells.Item(24, 20).Value = 2.25
$ws.Cells.Item(24, 21).Value = 1.85
$ws.Cells.Item(24, 22).Value = 1.95
$ws.Cells.Item(24, 23).Value = -1
$ws.Cells.Item(24, 24).Value = 2
$ws.Cells.Item(24, 25).Value = -1
$ws.Cells.Item(24, 26).Value = 0
$ws.Cells.Item(24, 27).Value = -0
$ws.Cells.Item(24, 28).Value = -1
$ws.Cells.Item(24, 29).Value = 0.95

# Row 25
$ws.Cells.Item(25, 2).Value = 6204317
$ws.Cells.Item(25, 6).Value = 'Jeonnam Dragons'
$ws.Cells.Item(25, 7).Value = 'Seoul ELand FC'
$ws.Cells.Item(25, 8).Value = 3
$ws.Cells.Item(25, 9).Value = 3
$ws.Cells.Item(25, 10).Value = 'D'
$ws.Cells.Item(25, 11).Value = 2.4
$ws.Cells.Item(25, 12).Value = 3.3
$ws.Cells.Item(25, 13).Value = 2.625
$ws.Cells.Item(25, 14).Value = 2.15
$ws.Cells.Item(25, 15).Value = 3.4
$ws.Cells.Item(25, 16).Value = 3
$ws.Cells.Item(25, 17).Value = -0.25
$ws.Cells.Item(25, 18).Value = 1.9
$ws.Cells.Item(25, 19).Value = 1.9
$ws.Cells.Item(25, 20).Value = 2.5
$ws.Cells.Item(25, 21).Value = 1.975
$ws.Cells.Item(25, 22).Value = 1.825
$ws.Cells.Item(25, 23).Value = -1
$ws.Cells.Item(25, 24).Value = 2.4
$ws.Cells.Item(25, 25).Value = -1
$ws.Cells.Item(25, 26).Value = -0.5
$ws.Cells.Item(25, 27).Value = 0.45
$ws.Cells.Item(25, 28).Value = 0.9750000000000001
$ws.Cells.Item(25, 29).Value = -1

# Row 74
$ws.Cells.Item(74, 2).Value = 6206241
$ws.Cells.Item(74, 6).Value = 'Gyeongnam FC'
$ws.Cells.Item(74, 7).Value = 'Seongnam FC'
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 2
$ws.Cells.Item(74, 10).Value = 'A'
$ws.Cells.Item(74, 11).Value = 2.3
$ws.Cells.Item(74, 12).Value = 3.4
$ws.Cells.Item(74, 13).Value = 2.75
$ws.Cells.Item(74, 14).Value = 2.15
$ws.Cells.Item(74, 15).Value = 3.4
$ws.Cells.Item(74, 16).Value = 3
$ws.Cells.Item(74, 17).Value = -0.25
$ws.Cells.Item(74, 18).Value = 1.85
$ws.Cells.Item(74, 19).Value = 1.95
$ws.Cells.Item(74, 20).Value = 2.5
$ws.Cells.Item(74, 21).Value = 1.8
$ws.Cells.Item(74, 22).Value = 2
$ws.Cells.Item(74, 23).Value = -1
$ws.Cells.Item(74, 24).Value = -1
$ws.Cells.Item(74, 25).Value = 2
$ws.Cells.Item(74, 26).Value = -1
$ws.Cells.Item(74, 27).Value = 0.95
$ws.Cells.Item(74, 28).Value = -1
$ws.Cells.Item(74, 29).Value = 1

# Row 75
$ws.Cells.Item(75, 2).Value = 6206240
$ws.Cells.Item(75, 6).Value = 'Gimpo FC'
$ws.Cells.Item(75, 7).Value = 'Busan I Park'
$ws.Cells.Item(75, 8).Value = 2
$ws.Cells.Item(75, 9).Value = 3
$ws.Cells.Item(75, 10).Value = 'A'
$ws.Cells.Item(75, 11).Value = 3.1
$ws.Cells.Item(75, 12).Value = 3.1
$ws.Cells.Item(75, 13).Value = 2.2
$ws.Cells.Item(75, 14).Value = 3.2
$ws.Cells.Item(75, 15).Value = 3
$ws.Cells.Item(75, 16).Value = 2.2
$ws.Cells.Item(75, 17).Value = 0.25
$ws.Cells.Item(75, 18).Value = 1.875
$ws.Cells.Item(75, 19).Value = 1.925
$ws.Cells.Item(75, 20).Value = 2
$ws.Cells.Item(75, 21).Value = 1.95
$ws.Cells.Item(75, 22).Value = 1.85
$ws.Cells.Item(75, 23).Value = -1
$ws.Cells.Item(75, 24).Value = -1
$ws.Cells.Item(75, 25).Value = 1.2
$ws.Cells.Item(75, 26).Value = -1
$ws.Cells.Item(75, 27).Value = 0.925
$ws.Cells.Item(75, 28).Value = 0.95
$ws.Cells.Item(75, 29).Value = -1

# Row 117
$ws.Cells.Item(117, 2).Value = 6352816
$ws.Cells.Item(117, 6).Value = 'Jeonnam Dragons'
$ws.Cells.Item(117, 7).Value = 'Ansan Greeners FC'
$ws.Cells.Item(117, 8).Value = 3
$ws.Cells.Item(117, 9).Value = 2
$ws.Cells.Item(117, 10).Value = 'H'
$ws.Cells.Item(117, 11).Value = 1.666
$ws.Cells.Item(117, 12).Value = 3.8
$ws.Cells.Item(117, 13).Value = 5
$ws.Cells.Item(117, 14).Value = 1.727
$ws.Cells.Item(117, 15).Value = 3.75
$ws.Cells.Item(117, 16).Value = 4.75
$ws.Cells.Item(117, 17).Value = -0.75
$ws.Cells.Item(117, 18).Value = 1.925
$ws.Cells.Item(117, 19).Value = 1.875
$ws.Cells.Item(117, 20).Value = 2.75
$ws.Cells.Item(117, 21).Value = 1.875
$ws.Cells.Item(117, 22).Value = 1.925
$ws.Cells.Item(117, 23).Value = 0.7270000000000001
$ws.Cells.Item(117, 24).Value = -1
$ws.Cells.Item(117, 25).Value = -1
$ws.Cells.Item(117, 26).Value = 0.4625
$ws.Cells.Item(117, 27).Value = -0.5
$ws.Cells.Item(117, 28).Value = 0.875
$ws.Cells.Item(117, 29).Value = -1

# Row 118
$ws.Cells.Item(118, 2).Value = 6414604
$ws.Cells.Item(118, 6).Value = 'Seongnam FC'
$ws.Cells.Item(118, 7).Value = 'Chungnam Asan FC'
$ws.Cells.Item(118, 8).Value = 2
$ws.Cells.Item(118, 9).Value = 0
$ws.Cells.Item(118, 10).Value = 'H'
$ws.Cells.Item(118, 11).Value = 2.4
$ws.Cells.Item(118, 12).Value = 3.2
$ws.Cells.Item(118, 13).Value = 3
$ws.Cells.Item(118, 14).Value = 2.4
$ws.Cells.Item(118, 15).Value = 3.2
$ws.Cells.Item(118, 16).Value = 3
$ws.Cells.Item(118, 17).Value = -0.25
$ws.Cells.Item(118, 18).Value = 2.05
$ws.Cells.Item(118, 19).Value = 1.75
$ws.Cells.Item(118, 20).Value = 2
$ws.Cells.Item(118, 21).Value = 1.825
$ws.Cells.Item(118, 22).Value = 1.975
$ws.Cells.Item(118, 23).Value = 1.4
$ws.Cells.Item(118, 24).Value = -1
$ws.Cells.Item(118, 25).Value = -1
$ws.Cells.Item(118, 26).Value = 1.05
$ws.Cells.Item(118, 27).Value = -1
$ws.Cells.Item(118, 28).Value = 0
$ws.Cells.Item(118, 29).Value = -0

# Row 140
$ws.Cells.Item(140, 2).Value = 6531883
$ws.Cells.Item(140, 6).Value = 'Seongnam FC'
$ws.Cells.Item(140, 7).Value = 'Ansan Greeners FC'
$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 9).Value = 2
$ws.Cells.Item(140, 10).Value = 'A'
$ws.Cells.Item(140, 11).Value = 1.7
$ws.Cells.Item(140, 12).Value = 3.8
$ws.Cells.Item(140, 13).Value = 4.2
$ws.Cells.Item(140, 14).Value = 1.75
$ws.Cells.Item(140, 15).Value = 4
$ws.Cells.Item(140, 16).Value = 3.75
$ws.Cells.Item(140, 17).Value = -0.75
$ws.Cells.Item(140, 18).Value = 2
$ws.Cells.Item(140, 19).Value = 1.8
$ws.Cells.Item(140, 20).Value = 3.25
$ws.Cells.Item(140, 21).Value = 2
$ws.Cells.Item(140, 22).Value = 1.8
$ws.Cells.Item(140, 23).Value = -1
$ws.Cells.Item(140, 24).Value = -1
$ws.Cells.Item(140, 25).Value = 2.75
$ws.Cells.Item(140, 26).Value = -1
$ws.Cells.Item(140, 27).Value = 0.8
$ws.Cells.Item(140, 28).Value = -1
$ws.Cells.Item(140, 29).Value = 0.8

# Row 141
$ws.Cells.Item(141, 2).Value = 6540655
$ws.Cells.Item(141, 6).Value = 'Gimpo FC'
$ws.Cells.Item(141, 7).Value = 'Gyeongnam FC'
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 9).Value = 1
$ws.Cells.Item(141, 10).Value = 'A'
$ws.Cells.Item(141, 11).Value = 2.4
$ws.Cells.Item(141, 12).Value = 3.3
$ws.Cells.Item(141, 13).Value = 2.7
$ws.Cells.Item(141, 14).Value = 2.75
$ws.Cells.Item(141, 15).Value = 3.25
$ws.Cells.Item(141, 16).Value = 2.3
$ws.Cells.Item(141, 17).Value = 0.25
$ws.Cells.Item(141, 18).Value = 1.75
$ws.Cells.Item(141, 19).Value = 2.05
$ws.Cells.Item(141, 20).Value = 2.25
$ws.Cells.Item(141, 21).Value = 1.825
$ws.Cells.Item(141, 22).Value = 1.975
$ws.Cells.Item(141, 23).Value = -1
$ws.Cells.Item(141, 24).Value = -1
$ws.Cells.Item(141, 25).Value = 1.3
$ws.Cells.Item(141, 26).Value = -1
$ws.Cells.Item(141, 27).Value = 1.05
$ws.Cells.Item(141, 28).Value = -1
$ws.Cells.Item(141, 29).Value = 0.9750000000000001

# Row 142
$ws.Cells.Item(142, 2).Value = 6537916
$ws.Cells.Item(142, 6).Value = 'Gimcheon Sangmu FC'
$ws.Cells.Item(142, 7).Value = 'Seoul ELand FC'
$ws.Cells.Item(142, 8).Value = 1
$ws.Cells.Item(142, 9).Value = 0
$ws.Cells.Item(142, 10).Value = 'H'
$ws.Cells.Item(142, 11).Value = 1.3
$ws.Cells.Item(142, 12).Value = 5
$ws.Cells.Item(142, 13).Value = 7.5
$ws.Cells.Item(142, 14).Value = 1.222
$ws.Cells.Item(142, 15).Value = 5.5
$ws.Cells.Item(142, 16).Value = 10
$ws.Cells.Item(142, 17).Value = -1.75
$ws.Cells.Item(142, 18).Value = 1.85
$ws.Cells.Item(142, 19).Value = 1.95
$ws.Cells.Item(142, 20).Value = 3.25
$ws.Cells.Item(142, 21).Value = 1.9
$ws.Cells.Item(142, 22).Value = 1.9
$ws.Cells.Item(142, 23).Value = 0.222
$ws.Cells.Item(142, 24).Value = -1
$ws.Cells.Item(142, 25).Value = -1
$ws.Cells.Item(142, 26).Value = -1
$ws.Cells.Item(142, 27).Value = 0.95
$ws.Cells.Item(142, 28).Value = -1
$ws.Cells.Item(142, 29).Value = 0.8999999999999999

# Row 143
$ws.Cells.Item(143, 2).Value = 6537915
$ws.Cells.Item(143, 6).Value = 'Busan I Park'
$ws.Cells.Item(143, 7).Value = 'Chungbuk Cheongju'
$ws.Cells.Item(143, 8).Value = 1
$ws.Cells.Item(143, 9).Value = 1
$ws.Cells.Item(143, 10).Value = 'D'
$ws.Cells.Item(143, 11).Value = 1.533
$ws.Cells.Item(143, 12).Value = 4
$ws.Cells.Item(143, 13).Value = 5.25
$ws.Cells.Item(143, 14).Value = 1.444
$ws.Cells.Item(143, 15).Value = 4.2
$ws.Cells.Item(143, 16).Value = 6
$ws.Cells.Item(143, 17).Value = -1.25
$ws.Cells.Item(143, 18).Value = 1.975
$ws.Cells.Item(143, 19).Value = 1.825
$ws.Cells.Item(143, 20).Value = 2.5
$ws.Cells.Item(143, 21).Value = 1.825
$ws.Cells.Item(143, 22).Value = 1.975
$ws.Cells.Item(143, 23).Value = -1
$ws.Cells.Item(143, 24).Value = 3.2
$ws.Cells.Item(143, 25).Value = -1
$ws.Cells.Item(143, 26).Value = -1
$ws.Cells.Item(143, 27).Value = 0.825
$ws.Cells.Item(143, 28).Value = -1
$ws.Cells.Item(143, 29).Value = 0.9750000000000001

# Row 144
$ws.Cells.Item(144, 2).Value = 6531882
$ws.Cells.Item(144, 6).Value = 'FC Anyang'
$ws.Cells.Item(144, 7).Value = 'Cheonan City'
$ws.Cells.Item(144, 8).Value = 2
$ws.Cells.Item(144, 9).Value = 1
$ws.Cells.Item(144, 10).Value = 'H'
$ws.Cells.Item(144, 11).Value = 1.571
$ws.Cells.Item(144, 12).Value = 4
$ws.Cells.Item(144, 13).Value = 4.75
$ws.Cells.Item(144, 14).Value = 1.4
$ws.Cells.Item(144, 15).Value = 4.333
$ws.Cells.Item(144, 16).Value = 6
$ws.Cells.Item(144, 17).Value = -1.25
$ws.Cells.Item(144, 18).Value = 2
$ws.Cells.Item(144, 19).Value = 1.8
$ws.Cells.Item(144, 20).Value = 2.75
$ws.Cells.Item(144, 21).Value = 1.975
$ws.Cells.Item(144, 22).Value = 1.825
$ws.Cells.Item(144, 23).Value = 0.3999999999999999
$ws.Cells.Item(144, 24).Value = -1
$ws.Cells.Item(144, 25).Value = -1
$ws.Cells.Item(144, 26).Value = -0.5
$ws.Cells.Item(144, 27).Value = 0.4
$ws.Cells.Item(144, 28).Value = 0.4875
$ws.Cells.Item(144, 29).Value = -0.5

# Row 145
$ws.Cells.Item(145, 2).Value = 6527572
$ws.Cells.Item(145, 6).Value = 'Bucheon'
$ws.Cells.Item(145, 7).Value = 'Jeonnam Dragons'
$ws.Cells.Item(145, 8).Value = 4
$ws.Cells.Item(145, 9).Value = 1
$ws.Cells.Item(145, 10).Value = 'H'
$ws.Cells.Item(145, 11).Value = 2.3
$ws.Cells.Item(145, 12).Value = 3.3
$ws.Cells.Item(145, 13).Value = 2.8
$ws.Cells.Item(145, 14).Value = 2.25
$ws.Cells.Item(145, 15).Value = 3.3
$ws.Cells.Item(145, 16).Value = 2.9
$ws.Cells.Item(145, 17).Value = -0.25
$ws.Cells.Item(145, 18).Value = 1.975
$ws.Cells.Item(145, 19).Value = 1.825
$ws.Cells.Item(145, 20).Value = 2.5
$ws.Cells.Item(145, 21).Value = 1.975
$ws.Cells.Item(145, 22).Value = 1.825
$ws.Cells.Item(145, 23).Value = 1.25
$ws.Cells.Item(145, 24).Value = -1
$ws.Cells.Item(145, 25).Value = -1
$ws.Cells.Item(145, 26).Value = 0.9750000000000001
$ws.Cells.Item(145, 27).Value = -1
$ws.Cells.Item(145, 28).Value = 0.9750000000000001
$ws.Cells.Item(145, 29).Value = -1

# Row 163
$ws.Cells.Item(163, 2).Value = 7738683
$ws.Cells.Item(163, 6).Value = 'Chungbuk Cheongju'
$ws.Cells.Item(163, 7).Value = 'FC Anyang'
$ws.Cells.Item(163, 8).Value = 1
$ws.Cells.Item(163, 9).Value = 1
$ws.Cells.Item(163, 10).Value = 'D'
$ws.Cells.Item(163, 11).Value = 2.3
$ws.Cells.Item(163, 12).Value = 3.2
$ws.Cells.Item(163, 13).Value = 2.7
$ws.Cells.Item(163, 14).Value = 2.75
$ws.Cells.Item(163, 15).Value = 3.2
$ws.Cells.Item(163, 16).Value = 2.3
$ws.Cells.Item(163, 17).Value = 0.25
$ws.Cells.Item(163, 18).Value = 1.75
$ws.Cells.Item(163, 19).Value = 2.05
$ws.Cells.Item(163, 20).Value = 2.25
$ws.Cells.Item(163, 21).Value = 1.925
$ws.Cells.Item(163, 22).Value = 1.875
$ws.Cells.Item(163, 23).Value = -1
$ws.Cells.Item(163, 24).Value = 2.2
$ws.Cells.Item(163, 25).Value = -1
$ws.Cells.Item(163, 26).Value = 0.375
$ws.Cells.Item(163, 27).Value = -0.5
$ws.Cells.Item(163, 28).Value = -0.5
$ws.Cells.Item(163, 29).Value = 0.4375

# Row 164
$ws.Cells.Item(164, 2).Value = 7737346
$ws.Cells.Item(164, 6).Value = 'Busan I Park'
$ws.Cells.Item(164, 7).Value = 'Gimpo FC'
$ws.Cells.Item(164, 8).Value = 0
$ws.Cells.Item(164, 9).Value = 1
$ws.Cells.Item(164, 10).Value = 'A'
$ws.Cells.Item(164, 11).Value = 1.8
$ws.Cells.Item(164, 12).Value = 3.25
$ws.Cells.Item(164, 13).Value = 4
$ws.Cells.Item(164, 14).Value = 1.7
$ws.Cells.Item(164, 15).Value = 3.3
$ws.Cells.Item(164, 16).Value = 4.5
$ws.Cells.Item(164, 17).Value = -0.75
$ws.Cells.Item(164, 18).Value = 2
$ws.Cells.Item(164, 19).Value = 1.8
$ws.Cells.Item(164, 20).Value = 2.25
$ws.Cells.Item(164, 21).Value = 1.9
$ws.Cells.Item(164, 22).Value = 1.9
$ws.Cells.Item(164, 23).Value = -1
$ws.Cells.Item(164, 24).Value = -1
$ws.Cells.Item(164, 25).Value = 3.5
$ws.Cells.Item(164, 26).Value = -1
$ws.Cells.Item(164, 27).Value = 0.8
$ws.Cells.Item(164, 28).Value = -1
$ws.Cells.Item(164, 29).Value = 0.8999999999999999

# Row 173
$ws.Cells.Item(173, 2).Value = 7738660
$ws.Cells.Item(173, 6).Value = 'Seoul ELand FC'
$ws.Cells.Item(173, 7).Value = 'Gimpo FC'
$ws.Cells.Item(173, 8).Value = 1
$ws.Cells.Item(173, 9).Value = 1
$ws.Cells.Item(173, 10).Value = 'D'
$ws.Cells.Item(173, 11).Value = 2.25
$ws.Cells.Item(173, 12).Value = 3.25
$ws.Cells.Item(173, 13).Value = 2.9
$ws.Cells.Item(173, 14).Value = 1.909
$ws.Cells.Item(173, 15).Value = 3.3
$ws.Cells.Item(173, 16).Value = 3.6
$ws.Cells.Item(173, 17).Value = -0.5
$ws.Cells.Item(173, 18).Value = 1.975
$ws.Cells.Item(173, 19).Value = 1.825
$ws.Cells.Item(173, 20).Value = 2.25
$ws.Cells.Item(173, 21).Value = 1.95
$ws.Cells.Item(173, 22).Value = 1.85
$ws.Cells.Item(173, 23).Value = -1
$ws.Cells.Item(173, 24).Value = 2.3
$ws.Cells.Item(173, 25).Value = -1
$ws.Cells.Item(173, 26).Value = -1
$ws.Cells.Item(173, 27).Value = 0.825
$ws.Cells.Item(173, 28).Value = -0.5
$ws.Cells.Item(173, 29).Value = 0.425

# Row 174
$ws.Cells.Item(174, 2).Value = 7737365
$ws.Cells.Item(174, 6).Value = 'Chungnam Asan FC'
$ws.Cells.Item(174, 7).Value = 'Seongnam FC'
$ws.Cells.Item(174, 8).Value = 1
$ws.Cells.Item(174, 9).Value = 1
$ws.Cells.Item(174, 10).Value = 'D'
$ws.Cells.Item(174, 11).Value = 2.2
$ws.Cells.Item(174, 12).Value = 3.25
$ws.Cells.Item(174, 13).Value = 3
$ws.Cells.Item(174, 14).Value = 2.1
$ws.Cells.Item(174, 15).Value = 3.25
$ws.Cells.Item(174, 16).Value = 3.25
$ws.Cells.Item(174, 17).Value = -0.25
$ws.Cells.Item(174, 18).Value = 1.85
$ws.Cells.Item(174, 19).Value = 1.95
$ws.Cells.Item(174, 20).Value = 2.5
$ws.Cells.Item(174, 21).Value = 1.95
$ws.Cells.Item(174, 22).Value = 1.85
$ws.Cells.Item(174, 23).Value = -1
$ws.Cells.Item(174, 24).Value = 2.25
$ws.Cells.Item(174, 25).Value = -1
$ws.Cells.Item(174, 26).Value = -0.5
$ws.Cells.Item(174, 27).Value = 0.475
$ws.Cells.Item(174, 28).Value = -1
$ws.Cells.Item(174, 29).Value = 0.8500000000000001

# Row 175
$ws.Cells.Item(175, 2).Value = 7738661
$ws.Cells.Item(175, 6).Value = 'Cheonan City'
$ws.Cells.Item(175, 7).Value = 'Gyeongnam FC'
$ws.Cells.Item(175, 8).Value = 2
$ws.Cells.Item(175, 9).Value = 2
$ws.Cells.Item(175, 10).Value = 'D'
$ws.Cells.Item(175, 11).Value = 4
$ws.Cells.Item(175, 12).Value = 3.4
$ws.Cells.Item(175, 13).Value = 1.833
$ws.Cells.Item(175, 14).Value = 3.6
$ws.Cells.Item(175, 15).Value = 3.3
$ws.Cells.Item(175, 16).Value = 1.95
$ws.Cells.Item(175, 17).Value = 0.5
$ws.Cells.Item(175, 18).Value = 1.8
$ws.Cells.Item(175, 19).Value = 2
$ws.Cells.Item(175, 20).Value = 2.5
$ws.Cells.Item(175, 21).Value = 2.025
$ws.Cells.Item(175, 22).Value = 1.775
$ws.Cells.Item(175, 23).Value = -1
$ws.Cells.Item(175, 24).Value = 2.3
$ws.Cells.Item(175, 25).Value = -1
$ws.Cells.Item(175, 26).Value = 0.8
$ws.Cells.Item(175, 27).Value = -1
$ws.Cells.Item(175, 28).Value = 1.025
$ws.Cells.Item(175, 29).Value = -1

# Row 176
$ws.Cells.Item(176, 2).Value = 7737347
$ws.Cells.Item(176, 6).Value = 'Busan I Park'
$ws.Cells.Item(176, 7).Value = 'Jeonnam Dragons'
$ws.Cells.Item(176, 8).Value = 0
$ws.Cells.Item(176, 9).Value = 1
$ws.Cells.Item(176, 10).Value = 'A'
$ws.Cells.Item(176, 11).Value = 1.909
$ws.Cells.Item(176, 12).Value = 3.25
$ws.Cells.Item(176, 13).Value = 3.75
$ws.Cells.Item(176, 14).Value = 1.909
$ws.Cells.Item(176, 15).Value = 3.25
$ws.Cells.Item(176, 16).Value = 3.8
$ws.Cells.Item(176, 17).Value = -0.5
$ws.Cells.Item(176, 18).Value = 1.925
$ws.Cells.Item(176, 19).Value = 1.875
$ws.Cells.Item(176, 20).Value = 2.25
$ws.Cells.Item(176, 21).Value = 1.95
$ws.Cells.Item(176, 22).Value = 1.85
$ws.Cells.Item(176, 23).Value = -1
$ws.Cells.Item(176, 24).Value = -1
$ws.Cells.Item(176, 25).Value = 2.8
$ws.Cells.Item(176, 26).Value = -1
$ws.Cells.Item(176, 27).Value = 0.875
$ws.Cells.Item(176, 28).Value = -1
$ws.Cells.Item(176, 29).Value = 0.8500000000000001

# Row 182
$ws.Cells.Item(182, 2).Value = 7737348
$ws.Cells.Item(182, 6).Value = 'Cheonan City'
$ws.Cells.Item(182, 7).Value = 'Busan I Park'
$ws.Cells.Item(182, 8).Value = 2
$ws.Cells.Item(182, 9).Value = 4
$ws.Cells.Item(182, 10).Value = 'A'
$ws.Cells.Item(182, 11).Value = 3
$ws.Cells.Item(182, 12).Value = 3.2
$ws.Cells.Item(182, 13).Value = 2.1
$ws.Cells.Item(182, 14).Value = 3.5
$ws.Cells.Item(182, 15).Value = 3.3
$ws.Cells.Item(182, 16).Value = 1.909
$ws.Cells.Item(182, 17).Value = 0.5
$ws.Cells.Item(182, 18).Value = 1.85
$ws.Cells.Item(182, 19).Value = 1.95
$ws.Cells.Item(182, 20).Value = 2.25
$ws.Cells.Item(182, 21).Value = 1.825
$ws.Cells.Item(182, 22).Value = 1.975
$ws.Cells.Item(182, 23).Value = -1
$ws.Cells.Item(182, 24).Value = -1
$ws.Cells.Item(182, 25).Value = 0.909
$ws.Cells.Item(182, 26).Value = -1
$ws.Cells.Item(182, 27).Value = 0.95
$ws.Cells.Item(182, 28).Value = 0.825
$ws.Cells.Item(182, 29).Value = -1

# Row 183
$ws.Cells.Item(183, 2).Value = 7738689
$ws.Cells.Item(183, 6).Value = 'Ansan Greeners FC'
$ws.Cells.Item(183, 7).Value = 'Chungnam Asan FC'
$ws.Cells.Item(183, 8).Value = 1
$ws.Cells.Item(183, 9).Value = 0
$ws.Cells.Item(183, 10).Value = 'H'
$ws.Cells.Item(183, 11).Value = 3
$ws.Cells.Item(183, 12).Value = 3
$ws.Cells.Item(183, 13).Value = 2.25
$ws.Cells.Item(183, 14).Value = 4
$ws.Cells.Item(183, 15).Value = 3.1
$ws.Cells.Item(183, 16).Value = 1.909
$ws.Cells.Item(183, 17).Value = 0.5
$ws.Cells.Item(183, 18).Value = 1.85
$ws.Cells.Item(183, 19).Value = 1.95
$ws.Cells.Item(183, 20).Value = 2.25
$ws.Cells.Item(183, 21).Value = 1.975
$ws.Cells.Item(183, 22).Value = 1.825
$ws.Cells.Item(183, 23).Value = 3
$ws.Cells.Item(183, 24).Value = -1
$ws.Cells.Item(183, 25).Value = -1
$ws.Cells.Item(183, 26).Value = 0.8500000000000001
$ws.Cells.Item(183, 27).Value = -1
$ws.Cells.Item(183, 28).Value = -1
$ws.Cells.Item(183, 29).Value = 0.825

# Remove now-unused trailing rows 190-195 (dataset now ends at row 189)
$ws.Rows("190:195").Delete()

Write-Output "applied"